# Updates the crypto price (D) and 1h volume change (E) columns to the
# latest scraped values (GitHub Actions cron refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.468.31'
$ws.Range("E2").Value = '  +1.94%  '
$ws.Range("D3").Value = '1.828.88'
$ws.Range("E3").Value = '  +1.93%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.31'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5138'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3907'
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07641'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("E10").Value = '  +0.94%  '
$ws.Range("E11").Value = '  +2.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.11'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.294'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.94%  '
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.553'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.00%  '
$ws.Range("D16").Value = '1.825.47'
$ws.Range("E16").Value = '  +1.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.45'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001082'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06670'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.72'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.191'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.18%  '
$ws.Range("D23").Value = '28.491.23'
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.16'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.256'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +7.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.80'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.66'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("D28").Value = '2.034.62'
$ws.Range("E28").Value = '  +1.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.397'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +4.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.14'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.77%  '
$ws.Range("E31").Value = '  +2.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1086'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.678'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.664'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07024'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2229'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.004'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +7.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02324'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.144'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6282'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.23'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.183'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.397'
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.48'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5905'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.710'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.13'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.983'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.199'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06922'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.70%  '
